# Auto-generated edit script applying the cryptos.xlsx price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores Price/Volume figures as literal text (e.g. "598.69",
# "27.00", "3.544.97") rather than numbers. A plain .Value assignment of a
# numeric-looking string auto-converts it to a real number (dropping
# trailing zeros / introducing float rounding), so for the cells whose new
# value parses as a plain number we force Text format first, assign the
# value, then restore the default "Normal" style so no stray number format
# is left on the cell.
$forceTextCells = @(
    "D5", "D6", "D10", "D11", "D14", "D16", "D19", "D20",
    "D21", "D22", "D25", "D28", "D31", "D33", "D34", "D35",
    "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45",
    "D46", "D47"
)
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.838.45"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.538.28"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "598.69"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").Value = "136.11"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "3.539.37"
$ws.Range("E7").Value = "  +3.22%  "
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").Value = "6.91"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("D13").Value = "4.144.82"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "0.0000181"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.544.97"
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "27.00"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "65.033.92"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("D20").Value = "14.32"
$ws.Range("E20").Value = "  +5.47%  "
$ws.Range("D21").Value = "5.81"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "386.67"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("E23").Value = "  +5.86%  "
$ws.Range("D24").Value = "3.687.72"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "74.24"
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +10.02%  "
$ws.Range("D28").Value = "7.61"
$ws.Range("E28").Value = "  +6.27%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +4.73%  "
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("D32").Value = "3.552.18"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").Value = "1.43"
$ws.Range("E33").Value = "  +21.31%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "23.86"
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("D37").Value = "169.80"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").Value = "6.92"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("E40").Value = "  +7.39%  "
$ws.Range("D41").Value = "0.0803"
$ws.Range("E41").Value = "  +5.85%  "
$ws.Range("D42").Value = "0.825"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("D43").Value = "26.76"
$ws.Range("E43").Value = "  +18.49%  "
$ws.Range("D44").Value = "42.61"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "4.44"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  +7.94%  "
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("E49").Value = "  +6.10%  "
$ws.Range("D50").Value = "2.443.96"
$ws.Range("E50").Value = "  +11.28%  "
$ws.Range("E51").Value = "  +13.25%  "

foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
